# Modify method to calculate eac of track.
#
# - "inf" sheet: the old crf_track / "Capital recovery factor of the track."
#   pair is replaced by interest_rate (value 0.08) plus a matching description;
#   the turnout-wage description gains a "(USD)" unit suffix; two new rows are
#   appended for the high-quality-track design parameters.
# - "deriv" sheet: a new max_path_difference parameter row is appended and the
#   sheet becomes the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "inf" sheet (sheetId 2)
# ---------------------------------------------------------------------------
$inf = $wb.Worksheets.Item("inf")

# Pass 1: names/values for the new rows and the renamed row
$inf.Range("A12").Value = "gross_tk_in_hq_track_lifetime"
$inf.Range("B12").Value = 200000000
$inf.Range("B12").NumberFormat = "#,##0"
$inf.Range("B12").HorizontalAlignment = -4108

$inf.Range("A13").Value = "high_quality_track_price"
$inf.Range("B13").Value = 800000
$inf.Range("B13").NumberFormat = "#,##0"
$inf.Range("B13").HorizontalAlignment = -4108

# Row 9: crf_track -> interest_rate
$inf.Range("A9").Value = "interest_rate"
$inf.Range("B9").Value = 0.08

# Pass 2: descriptions
$inf.Range("C9").Value = "Interest rate used to calculate capital recovery factor (rate)."
$inf.Range("C11").Value = "Wage cost of maintaining a turnout (USD)."
$inf.Range("C12").Value = "Design tons for high quality track. Gross tons that a hq track is suposed to support during its lifetime (gross ton-km)."
$inf.Range("C13").Value = "The price of 1km of hight quality track (USD/km)."

# Column B widens to fit the new, longer numeric entries
$inf.Columns.Item(2).ColumnWidth = 11.1

# Selection cursor follows the new last row
$inf.Range("A13").Select() | Out-Null

# ---------------------------------------------------------------------------
# "deriv" sheet (sheetId 3)
# ---------------------------------------------------------------------------
$deriv = $wb.Worksheets.Item("deriv")

# Row 7 (new): max_path_difference
$deriv.Range("A7").Value = "max_path_difference"
$deriv.Range("B7").Value = 0.5
$deriv.Range("B7").HorizontalAlignment = -4108
$deriv.Range("C7").Value = "Maximum difference in paths distance between rail and road options (coeff). Derivation won't happen if rail path is much longer than road path."

# "deriv" becomes the active / selected sheet, with C7 selected
$deriv.Activate()
$deriv.Range("C7").Select() | Out-Null
